# Update Leave Card 6/22/2023 5:35 PM
# Fills in the PERIOD (column A) dates and EARNED (column C) values for
# Table15 rows 59-90 on the "2018 LEAVE CREDITS" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Row -> date-serial (1st of each month, continuing the existing monthly
# PERIOD sequence already present through row 58 / 2023-02-01).
$periodDates = @{
    59 = 44986  # 2023-03-01
    60 = 45017  # 2023-04-01
    61 = 45047  # 2023-05-01
    62 = 45078  # 2023-06-01
    63 = 45108  # 2023-07-01
    64 = 45139  # 2023-08-01
    65 = 45170  # 2023-09-01
    66 = 45200  # 2023-10-01
    67 = 45231  # 2023-11-01
    68 = 45261  # 2023-12-01
    69 = 45292  # 2024-01-01
    70 = 45323  # 2024-02-01
    71 = 45352  # 2024-03-01
    72 = 45383  # 2024-04-01
    73 = 45413  # 2024-05-01
    74 = 45444  # 2024-06-01
    75 = 45474  # 2024-07-01
    76 = 45505  # 2024-08-01
    77 = 45536  # 2024-09-01
    78 = 45566  # 2024-10-01
    79 = 45597  # 2024-11-01
    80 = 45627  # 2024-12-01
    81 = 45658  # 2025-01-01
    82 = 45689  # 2025-02-01
    83 = 45717  # 2025-03-01
    84 = 45748  # 2025-04-01
    85 = 45778  # 2025-05-01
    86 = 45809  # 2025-06-01
    87 = 45839  # 2025-07-01
    88 = 45870  # 2025-08-01
    89 = 45901  # 2025-09-01
    90 = 45931  # 2025-10-01
}

foreach ($r in $periodDates.Keys) {
    $ws.Cells.Item($r, 1).Value = $periodDates[$r]
}

# The first three newly-dated rows (59-61) also have a recorded EARNED
# (VACATION LEAVE) value of 1.25 days, same as the prior rows in the table.
$ws.Cells.Item(59, 3).Value = 1.25
$ws.Cells.Item(60, 3).Value = 1.25
$ws.Cells.Item(61, 3).Value = 1.25
